$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Generic / consider using the original paragraph style" as DONE
$ws.Range("D27").Value = "DONE"
$ws.Rows(27).Hidden = $true

# Insert a new backlog row for "Scala CPS Plugin" right after the
# "Scala / partial functions" row, pushing the remaining rows down
$ws.Rows(29).Insert()
$ws.Range("A29").Value = "Scala CPS Plugin"
$ws.Range("B29").Value = "Reference programming with shift/reset"
$ws.Range("C29").Value = 15

# Mark "Scala.React / example: consistency" as DONE
$ws.Range("D30").Value = "DONE"
$ws.Rows(30).Hidden = $true

# Mark "Scala.React / example: lazy/strict signal, dependency" as DONE
$ws.Range("D31").Value = "DONE"
$ws.Rows(31).Hidden = $true

# Mark "Generic / change the UMLs to Visio PDFs" as DONE
$ws.Range("D33").Value = "DONE"
$ws.Rows(33).Hidden = $true

# Grow the table so it covers the newly inserted row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D34"))

# Leave the selection where the user ended up editing
$ws.Range("B34").Select()
